# Updates the cryptos list (prices, 1h volume %, and two swapped coin rows)
# as published by the "Updated cryptos list" GitHub Actions workflow run on
# Fri Jun 28 09:32:31 UTC 2024.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text value.
$updates = @(
    @{ Cell = "D2"; Value = "61.404.49" }
    @{ Cell = "E2"; Value = "  +0.67%  " }
    @{ Cell = "D3"; Value = "3.434.46" }
    @{ Cell = "E3"; Value = "  +1.38%  " }
    @{ Cell = "E4"; Value = "  +0.02%  " }
    @{ Cell = "D5"; Value = "575.51" }
    @{ Cell = "E5"; Value = "  +0.40%  " }
    @{ Cell = "D6"; Value = "145.07" }
    @{ Cell = "E6"; Value = "  +5.74%  " }
    @{ Cell = "D7"; Value = "3.434.69" }
    @{ Cell = "E7"; Value = "  +1.50%  " }
    @{ Cell = "E8"; Value = "  +0.03%  " }
    @{ Cell = "E9"; Value = "  +1.57%  " }
    @{ Cell = "E10"; Value = "  +0.36%  " }
    @{ Cell = "E11"; Value = "  +2.84%  " }
    @{ Cell = "E12"; Value = "  +1.47%  " }
    @{ Cell = "D13"; Value = "4.020.34" }
    @{ Cell = "E13"; Value = "  +1.45%  " }
    @{ Cell = "D14"; Value = "28.14" }
    @{ Cell = "E14"; Value = "  +5.92%  " }
    @{ Cell = "E15"; Value = "  -0.61%  " }
    @{ Cell = "E16"; Value = "  +0.94%  " }
    @{ Cell = "D17"; Value = "3.436.38" }
    @{ Cell = "E17"; Value = "  +1.55%  " }
    @{ Cell = "D18"; Value = "61.508.20" }
    @{ Cell = "E18"; Value = "  +0.77%  " }
    @{ Cell = "D19"; Value = "6.29" }
    @{ Cell = "E19"; Value = "  +7.51%  " }
    @{ Cell = "D20"; Value = "14.23" }
    @{ Cell = "E20"; Value = "  +2.61%  " }
    @{ Cell = "D21"; Value = "9.40" }
    @{ Cell = "E21"; Value = "  +1.09%  " }
    @{ Cell = "D22"; Value = "395.44" }
    @{ Cell = "E22"; Value = "  +5.47%  " }
    @{ Cell = "E23"; Value = "  +2.81%  " }
    @{ Cell = "D24"; Value = "73.72" }
    @{ Cell = "E24"; Value = "  +4.12%  " }
    @{ Cell = "D25"; Value = "0.996" }
    @{ Cell = "E25"; Value = "  -0.38%  " }
    @{ Cell = "E26"; Value = "  -0.28%  " }
    @{ Cell = "D27"; Value = "0.0000123" }
    @{ Cell = "E27"; Value = "  +0.24%  " }
    @{ Cell = "D28"; Value = "3.574.15" }
    @{ Cell = "E28"; Value = "  +1.80%  " }
    @{ Cell = "D29"; Value = "0.181" }
    @{ Cell = "E29"; Value = "  +4.36%  " }
    @{ Cell = "D30"; Value = "7.59" }
    @{ Cell = "E30"; Value = "  +3.09%  " }
    @{ Cell = "E31"; Value = "  +0.14%  " }
    @{ Cell = "D32"; Value = "8.25" }
    @{ Cell = "E32"; Value = "  +2.31%  " }
    @{ Cell = "E33"; Value = "  -8.90%  " }
    @{ Cell = "E34"; Value = "  +1.77%  " }
    @{ Cell = "E35"; Value = "  -0.03%  " }
    @{ Cell = "D36"; Value = "23.90" }
    @{ Cell = "E36"; Value = "  +2.23%  " }
    @{ Cell = "D37"; Value = "3.463.37" }
    @{ Cell = "E37"; Value = "  +1.71%  " }
    @{ Cell = "E38"; Value = "  +2.80%  " }
    @{ Cell = "B39"; Value = "ImmutableX" }
    @{ Cell = "C39"; Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx" }
    @{ Cell = "D39"; Value = "1.55" }
    @{ Cell = "E39"; Value = "  +0.10%  " }
    @{ Cell = "B40"; Value = "NEARProtocol" }
    @{ Cell = "C40"; Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near" }
    @{ Cell = "D40"; Value = "5.10" }
    @{ Cell = "E40"; Value = "  -0.25%  " }
    @{ Cell = "D41"; Value = "167.43" }
    @{ Cell = "E41"; Value = "  +1.62%  " }
    @{ Cell = "E42"; Value = "  +2.35%  " }
    @{ Cell = "D43"; Value = "27.05" }
    @{ Cell = "E43"; Value = "  +4.39%  " }
    @{ Cell = "E44"; Value = "  +3.30%  " }
    @{ Cell = "B45"; Value = "FirstDigitalUSD" }
    @{ Cell = "C45"; Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd" }
    @{ Cell = "D45"; Value = "1.00" }
    @{ Cell = "E45"; Value = "  +0.04%  " }
    @{ Cell = "B46"; Value = "Stacks" }
    @{ Cell = "C46"; Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx" }
    @{ Cell = "D46"; Value = "1.73" }
    @{ Cell = "E46"; Value = "  -0.46%  " }
    @{ Cell = "E47"; Value = "  +2.79%  " }
    @{ Cell = "D48"; Value = "42.26" }
    @{ Cell = "E48"; Value = "  +0.83%  " }
    @{ Cell = "D49"; Value = "2.598.00" }
    @{ Cell = "E49"; Value = "  +3.64%  " }
    @{ Cell = "E50"; Value = "  -1.39%  " }
    @{ Cell = "E51"; Value = "  +2.32%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Columns D (Price) and E (Volume) hold plain-text numeric-looking
    # strings (e.g. "6.29", dotted thousands like "61.404.49", or padded
    # percentages). Force text format before assignment so Excel keeps them
    # as text instead of auto-converting to numbers, then drop back to the
    # default "Normal" style so no stray formatting/style gets attached to
    # the cell (columns B/C - Coin/Link - do not need this, plain text is
    # safe there).
    if ($u.Cell -match "^[DE]") {
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
